# Auto-generated Excel COM-interop script applying the cell-value edits
# described by the upstream OOXML diff (scheduled profit-tracker refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2155.1
$ws.Range("I32").Value = 1766
$ws.Range("J32").Value = 2544.2
$ws.Range("K32").Value = 1766
$ws.Range("L32").Value = 2544.2
$ws.Range("M32").Value = -1440
$ws.Range("N32").Value = -3196.2

# Row 62
$ws.Range("H62").Value = 9762.182000000001
$ws.Range("I62").Value = 14749.167
$ws.Range("K62").Value = 14749.167
$ws.Range("M62").Value = -14125.167

# Row 65
$ws.Range("H65").Value = 9762.182000000001
$ws.Range("I65").Value = 14749.167
$ws.Range("K65").Value = 73745.83499999999
$ws.Range("M65").Value = -70625.83499999999

# Row 69
$ws.Range("H69").Value = 15999.4
$ws.Range("I69").Value = 9998.5
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 29995.5
$ws.Range("L69").Value = 60000
$ws.Range("M69").Value = -29121.5
$ws.Range("N69").Value = -61748

# Row 72
$ws.Range("H72").Value = 15999.4
$ws.Range("I72").Value = 9998.5
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 89986.5
$ws.Range("L72").Value = 180000
$ws.Range("M72").Value = -85618.5
$ws.Range("N72").Value = -188736

# Row 74
$ws.Range("H74").Value = 22097.334
$ws.Range("I74").Value = 6198.3335
$ws.Range("K74").Value = 6198.3335
$ws.Range("M74").Value = -5262.3335

# Row 77
$ws.Range("H77").Value = 22097.334
$ws.Range("I77").Value = 6198.3335
$ws.Range("K77").Value = 30991.6675
$ws.Range("M77").Value = -26311.6675

# Row 129
$ws.Range("H129").Value = 1102.5294
$ws.Range("I129").Value = 832.55554
$ws.Range("J129").Value = 1406.25
$ws.Range("K129").Value = 2497.66662
$ws.Range("L129").Value = 4218.75
$ws.Range("M129").Value = 2502.33338
$ws.Range("N129").Value = -14218.75

# Row 135
$ws.Range("H135").Value = 1070.3125
$ws.Range("I135").Value = 1008.1818
$ws.Range("K135").Value = 9073.636199999999
$ws.Range("M135").Value = -6538.636199999999

# Row 138
$ws.Range("H138").Value = 2579.6128
$ws.Range("J138").Value = 2716.5217
$ws.Range("L138").Value = 8149.5651
$ws.Range("N138").Value = -18429.5651

# Row 141
$ws.Range("H141").Value = 3502.6086
$ws.Range("I141").Value = 2049.3
$ws.Range("K141").Value = 6147.900000000001
$ws.Range("M141").Value = -967.9000000000005

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 4027.75
$ws.Range("I110").Value = 3703.6667
$ws.Range("K110").Value = 3703.6667
$ws.Range("M110").Value = -1658.6667

# Row 131
$ws.Range("H131").Value = 79999
$ws.Range("J131").Value = 79999
$ws.Range("L131").Value = 79999
$ws.Range("N131").Value = -90079

# Row 132
$ws.Range("H132").Value = 36806.875
$ws.Range("J132").Value = 6327.7393
$ws.Range("L132").Value = 18983.2179
$ws.Range("N132").Value = -24043.2179

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2481.6086
$ws.Range("I86").Value = 1405.3334
$ws.Range("J86").Value = 3173.5
$ws.Range("K86").Value = 1405.3334
$ws.Range("L86").Value = 3173.5
$ws.Range("M86").Value = -282.3334
$ws.Range("N86").Value = -5419.5

# Row 89
$ws.Range("H89").Value = 2481.6086
$ws.Range("I89").Value = 1405.3334
$ws.Range("J89").Value = 3173.5
$ws.Range("K89").Value = 7026.666999999999
$ws.Range("L89").Value = 15867.5
$ws.Range("M89").Value = -1410.666999999999
$ws.Range("N89").Value = -27099.5

# Row 108
$ws.Range("H108").Value = 80000
$ws.Range("J108").Value = 80000
$ws.Range("L108").Value = 80000
$ws.Range("N108").Value = -87680

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# Row 31
$ws.Range("H31").Value = 1420.15
$ws.Range("I31").Value = 1329.3214
$ws.Range("J31").Value = 1632.0834
$ws.Range("K31").Value = 1329.3214
$ws.Range("L31").Value = 1632.0834
$ws.Range("M31").Value = -1034.3214
$ws.Range("N31").Value = -2222.0834

# Row 34
$ws.Range("H34").Value = 1420.15
$ws.Range("I34").Value = 1329.3214
$ws.Range("J34").Value = 1632.0834
$ws.Range("K34").Value = 1329.3214
$ws.Range("L34").Value = 1632.0834
$ws.Range("M34").Value = -1127.3214
$ws.Range("N34").Value = -2036.0834

# Row 86
$ws.Range("H86").Value = 4871.2856
$ws.Range("I86").Value = 4849.8335
$ws.Range("K86").Value = 4849.8335
$ws.Range("M86").Value = -3726.8335

# Row 89
$ws.Range("H89").Value = 4871.2856
$ws.Range("I89").Value = 4849.8335
$ws.Range("K89").Value = 24249.1675
$ws.Range("M89").Value = -18633.1675

# Row 94
$ws.Range("H94").Value = 2701.7778
$ws.Range("I94").Value = 2650.75
$ws.Range("J94").Value = 2742.6
$ws.Range("K94").Value = 2650.75
$ws.Range("L94").Value = 2742.6
$ws.Range("M94").Value = -2199.75
$ws.Range("N94").Value = -3644.6

# Row 132
$ws.Range("H132").Value = 2225
$ws.Range("I132").Value = 2170
$ws.Range("K132").Value = 6510
$ws.Range("M132").Value = -3980

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 460.22223
$ws.Range("I121").Value = 465
$ws.Range("K121").Value = 1395
$ws.Range("M121").Value = -85

# Row 131
$ws.Range("H131").Value = 25036.23
$ws.Range("J131").Value = 24712.666
$ws.Range("L131").Value = 74137.99800000001
$ws.Range("N131").Value = -84217.99800000001

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3132.5
$ws.Range("I102").Value = 2875.5881
$ws.Range("K102").Value = 2875.5881
$ws.Range("M102").Value = -1253.5881

# Row 130
$ws.Range("H130").Value = 89998.5
$ws.Range("J130").Value = 89998.5
$ws.Range("L130").Value = 89998.5
$ws.Range("N130").Value = -100038.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4088.0833
$ws.Range("I40").Value = 3846.1365
$ws.Range("K40").Value = 3846.1365
$ws.Range("M40").Value = -3710.1365

# Row 69
$ws.Range("H69").Value = 40163
$ws.Range("J69").Value = 40163
$ws.Range("L69").Value = 40163
$ws.Range("N69").Value = -41785

# Row 72
$ws.Range("H72").Value = 40163
$ws.Range("J72").Value = 40163
$ws.Range("L72").Value = 120489
$ws.Range("N72").Value = -128601

# Row 122
$ws.Range("H122").Value = 4147
$ws.Range("I122").Value = 3933.875
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 11801.625
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -9351.625
$ws.Range("N122").Value = -19898.5

# Row 132
$ws.Range("H132").Value = 47449.332
$ws.Range("I132").Value = 76931.31
$ws.Range("J132").Value = 4566.4546
$ws.Range("K132").Value = 230793.93
$ws.Range("L132").Value = 13699.3638
$ws.Range("M132").Value = -228263.93
$ws.Range("N132").Value = -18759.3638

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1277.1
$ws.Range("I81").Value = 1277.1
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2554.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1493.2
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 1277.1
$ws.Range("I84").Value = 1277.1
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12771
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -7467
$ws.Range("N84").ClearContents()

# Row 108
$ws.Range("H108").Value = 48854.46
$ws.Range("J108").Value = 48854.46
$ws.Range("L108").Value = 48854.46
$ws.Range("N108").Value = -56534.46

# Row 119
$ws.Range("H119").Value = 27999.334
$ws.Range("J119").Value = 27999.334
$ws.Range("L119").Value = 27999.334
$ws.Range("N119").Value = -37675.334

# Row 133
$ws.Range("H133").Value = 81989.75
$ws.Range("J133").Value = 81989.75
$ws.Range("L133").Value = 81989.75
$ws.Range("N133").Value = -92109.75

